$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values before swapping (use Value2 to get actual values, not the Value property object)
$d5 = $ws.Range("D5").Value2
$d6 = $ws.Range("D6").Value2
$m6 = $ws.Range("M6").Value2
$u5 = $ws.Range("U5").Value2

# Swap item names between row 5 and row 6
$ws.Range("D5").Value = $d6
$ws.Range("D6").Value = $d5

# Row 5 now gets the quantity (3) previously in M6, placed in M5; U5 is cleared
$ws.Range("M5").Value = $m6
$ws.Range("U5").ClearContents()

# Row 6 now gets the quantity (344) previously in U5, placed in U6; M6 is cleared
$ws.Range("U6").Value = $u5
$ws.Range("M6").ClearContents()
